# Insert two new price-report rows (dated 2023-03-24 / serial 45009) above
# the current row 231, shifting all subsequent rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A231:A232").EntireRow.Insert()

# New row 231: Ajo Chino "1a (guarda)"
$ws.Range("A231").Value = 9
$ws.Range("B231").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C231").Value = "Metropolitana"
$ws.Range("D231").Value = 45009
$ws.Range("E231").Value = 13
$ws.Range("F231").Value = 100112003
$ws.Range("G231").Value = "Ajo"
$ws.Range("H231").Value = "Chino"
$ws.Range("I231").Value = "1a (guarda)"
$ws.Range("J231").Value = 160
$ws.Range("K231").Value = 16000
$ws.Range("L231").Value = 16000
$ws.Range("M231").Value = 16000
$ws.Range("N231").Value = "`$/caja 10 kilos"
$ws.Range("O231").Value = "Provincia de Talagante"
$ws.Range("P231").Value = 1600
$ws.Range("Q231").Value = 10
$ws.Range("R231").Value = "Hortaliza"

# New row 232: Ajo Chino "2a (guarda)"
$ws.Range("A232").Value = 9
$ws.Range("B232").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C232").Value = "Metropolitana"
$ws.Range("D232").Value = 45009
$ws.Range("E232").Value = 13
$ws.Range("F232").Value = 100112003
$ws.Range("G232").Value = "Ajo"
$ws.Range("H232").Value = "Chino"
$ws.Range("I232").Value = "2a (guarda)"
$ws.Range("J232").Value = 70
$ws.Range("K232").Value = 15000
$ws.Range("L232").Value = 15000
$ws.Range("M232").Value = 15000
$ws.Range("N232").Value = "`$/caja 10 kilos"
$ws.Range("O232").Value = "Provincia de Talagante"
$ws.Range("P232").Value = 1500
$ws.Range("Q232").Value = 10
$ws.Range("R232").Value = "Hortaliza"
